# Update the "想去人数" (F column) counts on the 展览 (sheet1) and
# 全部类型 (sheet4) worksheets to reflect the newly generated data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    6  = 34
    8  = 1769
    9  = 50
    12 = 2025
    13 = 19
    14 = 141
    15 = 1288
    16 = 461
    17 = 18
    18 = 289
    26 = 1095
    28 = 329
    31 = 317
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "全部类型" - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    6  = 34
    8  = 1769
    10 = 50
    13 = 2025
    14 = 19
    15 = 141
    16 = 1288
    17 = 461
    18 = 18
    19 = 289
    27 = 1095
    29 = 329
    32 = 317
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
